$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A122").Value = 'Objęto go ochroną w celach naukowo-dydaktycznych, ze względu na atrakcyjny pokrój drzewa oraz budowę morfologiczną liści i kory drzewa, a także ze względu na pokaźne rozmiary.'
$ws.Range("B122").Value = 'polish'
$ws.Range("A123").Value = 'Stan zdrowotny drzewa jest zadowalający, aczkolwiek nie jest ono zachowane w całości – pozostała głównie kłoda i ok. 40% masy korony, brak jest partii wierzchołkowej, spowodowany złamaniem przewodnika oraz wykonaniem cięć redukcyjnych.'
$ws.Range("B123").Value = 'polish'
$ws.Range("A124").Value = 'Drużynowy konkurs skoków narciarskich na normalnej skoczni powrócił po tym, jak nie odbył się on na Mistrzostwach Świata w 2003.'
$ws.Range("B124").Value = 'polish'
$ws.Range("A125").Value = 'Miał obszerną bibliotekę specjalistycznych książek, literatury pięknej, a także dzieł słowiańskich.'
$ws.Range("B125").Value = 'polish'
$ws.Range("A126").Value = 'W czasie wojny domowej miejscowość znalazła się w rejonie działania antyrządowych bojówek.'
$ws.Range("B126").Value = 'polish'
$ws.Range("A127").Value = 'W 2020 roku miejscowe władze odkryły pokaźny arsenał terrorystów, w tym karabiny automatyczne, snajperskie, pociski moździerzowe, rakiety, urządzenia wybuchowe i znaczną ilość amunicji'
$ws.Range("B127").Value = 'polish'
$ws.Range("A128").Value = 'Zgodnie z tradycją co roku w mieście organizowany jest festiwal ziemniaków.'
$ws.Range("B128").Value = 'polish'
$ws.Range("A129").Value = 'Astronautyka obejmuje poznawanie oraz analizowanie warunków i zjawisk towarzyszących lotom statków kosmicznych. '
$ws.Range("B129").Value = 'polish'
$ws.Range("A130").Value = 'Odkrywa możliwości techniczne ich realizacji oraz bada oddziaływanie warunków lotu na psychofizyczny stan astronautów i możliwość ich adaptacji do tych warunków.'
$ws.Range("B130").Value = 'polish'
$ws.Range("A131").Value = 'W rozwiązywaniu zagadnień z zakresu astronautyki korzysta się z osiągnięć wielu gałęzi nauki (aerodynamiki, magnetohydrodynamiki, aeronomii, astrofizyki, planetologii, techniki rakietowej i innych).'
$ws.Range("B131").Value = 'polish'
$ws.Range("A132").Value = 'Idea podróży kosmicznych nurtowała ludzkie umysły od stuleci, ale pozostawała marzeniem aż do czasu zbudowania potężnych rakiet zdolnych unieść użyteczny ładunek daleko w przestrzeń.'
$ws.Range("B132").Value = 'polish'
$ws.Range("A133").Value = 'W ciągu niespełna trzydziestu lat sondy kosmiczne przeleciały obok wszystkich planet Układu Słonecznego.'
$ws.Range("B133").Value = 'polish'
$ws.Range("A134").Value = 'Tańce mogą być podzielone i opisane ze względu na rodzaj choreografii, rodzaj ruchów lub historyczne pochodzenie.'
$ws.Range("B134").Value = 'polish'
$ws.Range("A135").Value = 'Przed wynalezieniem pisma taniec był ważnym czynnikiem międzypokoleniowego przekazu historycznego.'
$ws.Range("B135").Value = 'polish'
$ws.Range("A136").Value = 'Najprostszy podział, jakiego można dokonać, to wyodrębnienie trzech odmiennych sfer działalności ludzkiej, w których się rozwinął: religia, rozrywka i sztuka.'
$ws.Range("B136").Value = 'polish'
$ws.Range("A137").Value = 'Poczta w znaczeniu nowożytnym, będąca częścią europejskiej sieci pocztowej, powstała za czasów króla Zygmunta Augusta.'
$ws.Range("B137").Value = 'polish'
$ws.Range("A138").Value = 'Zaprezentowany dwa lata później Szantaż był pierwszym brytyjskim filmem dźwiękowym.'
$ws.Range("B138").Value = 'polish'
$ws.Range("A139").Value = 'Arboretum w Rogowie – należy do grupy najcenniejszych, najbogatszych w gatunki i odmiany drzew i krzewów tego typu ogrodów w Europie.'
$ws.Range("B139").Value = 'polish'
$ws.Range("A140").Value = 'Obok jeziora Żółkin przebiega XII trasa rowerowa Bełchatowskiego Stowarzyszenia Załogi Rowerowej "Zgrzyt".'
$ws.Range("B140").Value = 'polish'
$ws.Range("A141").Value = 'Nadleśnictwo położone jest w północno-zachodniej części terenu województwa wielkopolskiego na terenie 5 powiatów: obornickim, szamotulskim, czarnkowsko-trzcianeckim, chodzieskim i poznańskim.'
$ws.Range("B141").Value = 'polish'
$ws.Range("A142").Value = 'Macedonia was an ancient kingdom on the periphery of Archaic and Classical Greece, and later the dominant state of Hellenistic Greece.'
$ws.Range("B142").Value = 'english'
$ws.Range("A143").Value = 'Canada bans the sale of assault-style firearms after a mass killing in Nova Scotia leaves twenty-two victims dead.'
$ws.Range("B143").Value = 'english'
$ws.Range("A144").Value = 'The remoteness of the Glensanda settlement is such that there are no road, rail, or marked footway links across the granite mountain, moor, heather and peat bog of the private Glensanda estate.'
$ws.Range("B144").Value = 'english'
$ws.Range("A145").Value = 'The only practical access is by boat from the shores of Loch Linnhe.'
$ws.Range("B145").Value = 'english'
$ws.Range("A146").Value = ' The band was playing the club circuit in London, and he was invited to join them for a few songs.'
$ws.Range("B146").Value = 'english'
$ws.Range("A147").Value = 'There are roads and private residences along much of the lake, except the southeastern side where the Caribou Bog complex borders the lake.'
$ws.Range("B147").Value = 'english'
$ws.Range("A148").Value = 'The squadron returned to the United States in December 1945 and was inactivated in March 1946, and its personnel and equipment transferred to another organization.'
$ws.Range("B148").Value = 'english'
$ws.Range("A149").Value = 'Organdy or organdie is the sheerest and crispest cotton cloth made.'
$ws.Range("B149").Value = 'english'
$ws.Range("A150").Value = 'The latter two finishes are more popular for summer wear and draped apparel whereas the first is more popular for loose apparel and home textiles such as dresses and curtains.'
$ws.Range("B150").Value = 'english'
$ws.Range("A151").Value = 'The film also denounces gang violence and presents meaningful solutions from former gang-members to stop this problem.'
$ws.Range("B151").Value = 'english'
$ws.Range("A152").Value = 'Her collected writings and recollections of that period offer a rare first hand English language accounts of life in wartime Japan.'
$ws.Range("B152").Value = 'english'
$ws.Range("A153").Value = 'Zanskar, ringed by high Himalayan mountains in northwest India, one of the most remote places on the planet, has been safe until now.'
$ws.Range("B153").Value = 'english'
$ws.Range("A154").Value = 'While they complete the school they are also placing local children in other schools and monasteries in the city of Manali and beyond. '
$ws.Range("B154").Value = 'english'
$ws.Range("A155").Value = 'After this, the Serbian player pushed Shengelia and the Georgian forward reacted with a punch. '
$ws.Range("B155").Value = 'english'
$ws.Range("A156").Value = 'After this, all the players who were in the bench came into the court and started the brawl.'
$ws.Range("B156").Value = 'english'
$ws.Range("A157").Value = 'Waller plans on sponsoring legislation to require disclosure on clothes collection bins operated by for-profit companies, to include a new derivative of methamphetamine to Colorado''s list of banned drugs, and to exclude dry-ice bombs from a felony list of explosives, and to allow homeschooled students to enroll in college classes.'
$ws.Range("B157").Value = 'english'
$ws.Range("A158").Value = 'On his dismissal as governor, he said of the State treasury that he left it empty because he met it empty.'
$ws.Range("B158").Value = 'english'
$ws.Range("A159").Value = 'In June 2014, its assets and infrastructure were purchased by a new company.'
$ws.Range("B159").Value = 'english'
$ws.Range("A160").Value = 'The village has a population of 589.'
$ws.Range("B160").Value = 'english'
$ws.Range("A161").Value = 'That is, the information is intended to be quickly found when needed.'
$ws.Range("B161").Value = 'english'
$ws.Range("A162").Value = 'La biología se ocupa tanto de la descripción de las características y los comportamientos de los organismos individuales, como de las especies en su conjunto, así como de la reproducción de los seres vivos y de las interacciones entre ellos y el entorno.'
$ws.Range("B162").Value = 'spanish'
$ws.Range("A163").Value = 'La biología moderna se divide en sub-disciplinas según los tipos de organismos y la escala en que se los estudia. '
$ws.Range("B163").Value = 'spanish'
$ws.Range("A164").Value = 'El virus se transmite generalmente de persona a persona a través de las pequeñas gotas de saliva.'
$ws.Range("B164").Value = 'spanish'
$ws.Range("A165").Value = 'Se han cerrado colegios y universidades en más de 124 países, lo que ha afectado a más de 2200 millones de estudiantes.'
$ws.Range("B165").Value = 'spanish'
$ws.Range("A166").Value = 'La mayoría de las provincias implementaron cierres de escuelas y guarderías, prohibiciones de grandes reuniones, así como el cierre de varios lugares de ocio y entretenimiento a mediados de marzo.'
$ws.Range("B166").Value = 'spanish'
$ws.Range("A167").Value = 'En Corea del Norte no se ha registrado ningún caso oficial, lo cual fue muy cuestionado debido a sus fronteras con China y Corea del Sur.'
$ws.Range("B167").Value = 'spanish'
$ws.Range("A168").Value = 'Puedes avisarnos de un artículo de actualidad que falte.'
$ws.Range("B168").Value = 'spanish'
$ws.Range("A169").Value = 'Los componentes del movimiento, aunque conservan su independencia en los límites del estatuto del movimiento, actúan siempre de conformidad con sus principios fundamentales y colaboran entre sí en el desempeño de sus tareas respectivas y para realizar su misión común.'
$ws.Range("B169").Value = 'spanish'
$ws.Range("A170").Value = 'En la segunda mitad del siglo XIX, la rápida evolución de la tecnología de las armas de fuego ocasionó un dramático incremento del número de muertos y de heridos en tiempo de guerra.'
$ws.Range("B170").Value = 'spanish'
$ws.Range("A171").Value = 'El símbolo debía ser sencillo, identificable a distancia, conocido por todos e idéntico para amigos y adversarios. '
$ws.Range("B171").Value = 'spanish'
$ws.Range("A172").Value = 'En el primer Convenio de Ginebra, se eligió la cruz roja sobre fondo blanco como emblema distintivo único.  '
$ws.Range("B172").Value = 'spanish'
$ws.Range("A173").Value = ' Considera también revisar los espacios de nombres distintos al principal, especialmente los anexos, aquí.'
$ws.Range("B173").Value = 'spanish'
$ws.Range("A174").Value = 'En la ciénega se reúnen varios arroyos con los que dan vida al río Lauca.'
$ws.Range("B174").Value = 'spanish'
$ws.Range("A175").Value = 'Una vez calculados para cada mes, ambos valores son calculados para todo el año y pueden ser leídos en la columna vertical al lado derecho del diagrama.'
$ws.Range("B175").Value = 'spanish'
$ws.Range("A176").Value = 'Sus contribuciones en econometría incluir el estudio de riesgo y seguros en los países en desarrollo.'
$ws.Range("B176").Value = 'spanish'
$ws.Range("A177").Value = 'Tepito es una estación del Metro de Ciudad de México perteneciente a la línea B, localizada en el barrio de Tepito.'
$ws.Range("B177").Value = 'spanish'
$ws.Range("A178").Value = 'Causa pudrición de las raíces en muchas especies de plantas.'
$ws.Range("B178").Value = 'spanish'
$ws.Range("A179").Value = 'Carne blanquecina y firme en el sombrero, leñosa y fibrosa en el pie, de sabor suave en los ejemplares jóvenes, amargo y desagradable en los ejemplares adultos, y de olor fuerte no muy agradable.  '
$ws.Range("B179").Value = 'spanish'
$ws.Range("A180").Value = 'Por lo tanto los dos sustratos de esta enzima son nitrito, y iones hidrógeno; mientras que sus tres productos son óxido nítrico, nitrato, y agua.'
$ws.Range("B180").Value = 'spanish'
$ws.Range("A181").Value = 'El anfitrión provee nutrientes orgánicos al hongo, y la actividad metabólica dentro del complejo es considerablemente mayor que afuera.'
$ws.Range("B181").Value = 'spanish'

$null = $ws.Range("C126").Select()

